$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("F4").Value = 1.0371
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 0.6901
$ws.Range("I4").Value = 0

# Row 5
$ws.Range("F5").Value = 1.106
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0.8909
$ws.Range("I5").Value = 0

# Row 6
$ws.Range("C6").Value = 0.0254
$ws.Range("D6").Value = 76
$ws.Range("E6").Value = 0.3188
$ws.Range("F6").Value = 0.4747
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 59276.1261
$ws.Range("I6").Value = 0

# Row 7
$ws.Range("C7").Value = 0.0097
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.1391
$ws.Range("F7").Value = 0.4333
$ws.Range("H7").Value = 2.8263
$ws.Range("I7").Value = 0

# Row 8
$ws.Range("C8").Value = 0.0291
$ws.Range("D8").Value = 347
$ws.Range("E8").Value = 0.2006
$ws.Range("F8").Value = 0.3107
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 2627.2297
$ws.Range("I8").Value = 0

# Row 9
$ws.Range("F9").Value = 0.3198
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0.7712
$ws.Range("I9").Value = 1

# Row 10
$ws.Range("F10").Value = 0.4554
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.9021
$ws.Range("I10").Value = 1

# Row 11
$ws.Range("F11").Value = 0.1255
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0.5069
$ws.Range("I11").Value = 0

# Row 12
$ws.Range("F12").Value = 0.5141
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0.8048
$ws.Range("I12").Value = 1

# Row 13
$ws.Range("D13").Value = 5
$ws.Range("F13").Value = 0.4341
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 3.9187
$ws.Range("I13").Value = 0

# Row 14
$ws.Range("F14").Value = 0.4527
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0.4627
$ws.Range("I14").Value = 0

# Row 15
$ws.Range("F15").Value = 0.4989
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0.4639

# Row 16
$ws.Range("F16").Value = 0.2138
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 3.2094
$ws.Range("I16").Value = 0

# Row 17
$ws.Range("F17").Value = 0.3027
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 4.7485
$ws.Range("I17").Value = 1

# Row 18
$ws.Range("F18").Value = 0.3181
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 148069.2947
$ws.Range("I18").Value = 1

# Row 19
$ws.Range("F19").Value = 0.3489
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 4.6965
$ws.Range("I19").Value = 0

# Row 20
$ws.Range("F20").Value = 0.4086
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 154244.6426
$ws.Range("I20").Value = 0
